$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Instance (D2) from "Automation3" to "Automation2"
$ws.Range("D2").Value = "Automation2"

# Update the selected cell / active cell on the sheet from B3 to E7
$ws.Range("E7").Select()
